$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - Tuesday, Jan 10 departure (Wizz Air A21N to London)
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Tuesday, Jan 10"
$ws.Range("C11").Value = "10:10 AM"
$ws.Range("D11").Value = "W92065"
$ws.Range("E11").Value = "London"
$ws.Range("F11").Value = "(LTN)"
$ws.Range("G11").Value = "Wizz Air "
$ws.Range("H11").Value = "A21N"
$ws.Range("I11").Value = "(G-WUKM)"
$ws.Range("J11").Value = "10:31 AM"
$ws.Range("K11").Font.Bold = $false
$ws.Range("L11").Value = "0 hours, 21 minutes"
$ws.Range("M11").Font.Bold = $false

# Row 12 - Tuesday, Jan 10 departure (Wizz Air A321 to Oslo)
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Tuesday, Jan 10"
$ws.Range("C12").Value = "11:05 AM"
$ws.Range("D12").Value = "W62079"
$ws.Range("E12").Value = "Oslo"
$ws.Range("F12").Value = "(TRF)"
$ws.Range("G12").Value = "Wizz Air "
$ws.Range("H12").Value = "A321"
$ws.Range("I12").Value = "(HA-LTB)"
$ws.Range("J12").Value = "11:14 AM"
$ws.Range("K12").Font.Bold = $false
$ws.Range("L12").Value = "0 hours, 9 minutes"
$ws.Range("M12").Font.Bold = $false
